$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.045.51'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.301.73'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '300.75'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.53'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E8').ClearFormats()

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.73'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0793'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '49.14'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.75%  '
$ws.Range('E12').ClearFormats()

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.18'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +13.61%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.77'
$ws.Range('D15').ClearFormats()

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.656.12'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.298.85'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('E17').ClearFormats()

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.29%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '43.003.83'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.73'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('E20').ClearFormats()

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.06'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.62'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.74%  '
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '236.78'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.04'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +6.47%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '24.47'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.28'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +4.90%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '166.61'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '33.92'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('E31').ClearFormats()

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('E32').ClearFormats()

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('E33').ClearFormats()

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.64'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +7.65%  '
$ws.Range('E35').ClearFormats()

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('E36').ClearFormats()

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '16.75'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +4.49%  '
$ws.Range('E37').ClearFormats()

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.83'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('E39').ClearFormats()

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('E40').ClearFormats()

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('E41').ClearFormats()

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('E42').ClearFormats()

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.983.86'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.64%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0284'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.95'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +3.23%  '
$ws.Range('E46').ClearFormats()

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.85'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.534.23'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '53.22'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('E50').ClearFormats()

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.58%  '
$ws.Range('E51').ClearFormats()

